$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3. Well Builder")
$ws.Activate()

# Insert two new rows above row 8 to make room for the new 1.6/1.7 sub-items,
# pushing the former rows 8-10 down to 10-12.
$ws.Rows("8:9").Insert()

# New sub-item rows (1.6 / 1.7) under "Build main GUI window..." (section 1).
$ws.Range("C8").Value = 1.6
$ws.Range("D8").Value = "Remove Casing"

$ws.Range("C9").Value = 1.7

# New top-level items 5 and 6.
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "Submit button save into database "

$ws.Range("B14").Value = 6
$ws.Range("C14").Value = "Update casing info based on casing selected"

# Filled after the items above, matching the shared-string allocation order
# recorded in the saved workbook.
$ws.Range("D9").Value = "Change casing"

# Match the saved selection state.
$ws.Range("C22").Select()
